# The commit swaps the contents of ppt/theme/theme1.xml (the deck's
# main/SlideMaster theme, currently the "Integral" colour scheme) and
# ppt/theme/theme2.xml (the NotesMaster theme, currently the default
# "Office Theme" colour scheme) - i.e. theme1.xml ends up holding the
# Office Theme colours and theme2.xml ends up holding the Integral
# colours. Font scheme / format scheme (fills, lines, effects) are
# identical between the two themes, so the only substantive change is
# the 12-colour <a:clrScheme>.
#
# The PowerPoint object model exposes the document theme's colour
# scheme through Master.Theme.ThemeColorScheme, whose 12 slots map 1:1
# onto dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink in that order. We drive
# the SlideMaster's theme (the one that actually paints the slides) to
# the "Office Theme" palette that theme2.xml currently holds.
#
# ColorFormat.RGB takes a single packed integer (R + G*256 + B*65536),
# the same encoding VBA's RGB() builtin produces, spelled out in hex
# below as (B*0x10000 + G*0x100 + R) so each target colour stays easy
# to read/verify against the srgbClr hex values in the diff.

$p = $ppt.ActivePresentation
$tcs = $p.SlideMaster.Theme.ThemeColorScheme

$tcs.Colors(1).RGB  = 0x000000   # dk1      -> 000000
$tcs.Colors(2).RGB  = 0xFFFFFF   # lt1      -> FFFFFF
$tcs.Colors(3).RGB  = 0x6A5444   # dk2      -> 44546A
$tcs.Colors(4).RGB  = 0xE6E6E7   # lt2      -> E7E6E6
$tcs.Colors(5).RGB  = 0xD59B5B   # accent1  -> 5B9BD5
$tcs.Colors(6).RGB  = 0x317DED   # accent2  -> ED7D31
$tcs.Colors(7).RGB  = 0xA5A5A5   # accent3  -> A5A5A5
$tcs.Colors(8).RGB  = 0x00C0FF   # accent4  -> FFC000
$tcs.Colors(9).RGB  = 0xC47244   # accent5  -> 4472C4
$tcs.Colors(10).RGB = 0x47AD70   # accent6  -> 70AD47
$tcs.Colors(11).RGB = 0xC16305   # hlink    -> 0563C1
$tcs.Colors(12).RGB = 0x724F95   # folHlink -> 954F72
